# Auto-generated: apply updated market/profit values to the Chocobo_Profits workbook
# (scheduled runner refresh of currentAveragePrice* / Leve* computed columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 6131.8945
$ws.Range("I98").Value = 4133
$ws.Range("J98").Value = 7435.522
$ws.Range("K98").Value = 4133
$ws.Range("L98").Value = 7435.522
$ws.Range("M98").Value = -2635
$ws.Range("N98").Value = -10431.522
$ws.Range("H122").Value = 6131.8945
$ws.Range("I122").Value = 4133
$ws.Range("J122").Value = 7435.522
$ws.Range("K122").Value = 12399
$ws.Range("L122").Value = 22306.566
$ws.Range("M122").Value = -9949
$ws.Range("N122").Value = -27206.566
$ws.Range("H123").Value = 41137.273
$ws.Range("J123").Value = 41137.273
$ws.Range("L123").Value = 41137.273
$ws.Range("N123").Value = -50937.273
$ws.Range("H128").Value = 41868.332
$ws.Range("J128").Value = 41868.332
$ws.Range("L128").Value = 41868.332
$ws.Range("N128").Value = -51828.332
$ws.Range("H132").Value = 26002402
$ws.Range("I132").Value = 35859390
$ws.Range("K132").Value = 107578170
$ws.Range("M132").Value = -107575640
$ws.Range("H134").Value = 49333.21
$ws.Range("J134").Value = 49333.21
$ws.Range("L134").Value = 49333.21
$ws.Range("N134").Value = -59473.21
$ws.Range("H140").Value = 71673.78
$ws.Range("J140").Value = 71673.78
$ws.Range("L140").Value = 71673.78
$ws.Range("N140").Value = -82033.78
$ws.Range("H141").Value = 6126.2607
$ws.Range("I141").Value = 6352.5713
$ws.Range("J141").Value = 3750
$ws.Range("K141").Value = 19057.7139
$ws.Range("L141").Value = 11250
$ws.Range("M141").Value = -13877.7139
$ws.Range("N141").Value = -21610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 10000
$ws.Range("K3").Value = 10000
$ws.Range("M3").Value = -9885
$ws.Range("H32").Value = 4899.649
$ws.Range("I32").Value = 4720.8887
$ws.Range("K32").Value = 4720.8887
$ws.Range("M32").Value = -4433.8887
$ws.Range("H45").Value = 2048
$ws.Range("I45").Value = 2302.4285
$ws.Range("K45").Value = 2302.4285
$ws.Range("M45").Value = -1925.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 192153.31
$ws.Range("I31").Value = 423096.84
$ws.Range("J31").Value = 2661.1794
$ws.Range("K31").Value = 423096.84
$ws.Range("L31").Value = 2661.1794
$ws.Range("M31").Value = -422801.84
$ws.Range("N31").Value = -3251.1794
$ws.Range("H34").Value = 192153.31
$ws.Range("I34").Value = 423096.84
$ws.Range("J34").Value = 2661.1794
$ws.Range("K34").Value = 423096.84
$ws.Range("L34").Value = 2661.1794
$ws.Range("M34").Value = -422894.84
$ws.Range("N34").Value = -3065.1794
$ws.Range("H99").Value = 4216.4165
$ws.Range("I99").Value = 1751.4
$ws.Range("J99").Value = 5977.143
$ws.Range("K99").Value = 1751.4
$ws.Range("L99").Value = 5977.143
$ws.Range("M99").Value = -253.4000000000001
$ws.Range("N99").Value = -8973.143
$ws.Range("H126").Value = 4216.4165
$ws.Range("I126").Value = 1751.4
$ws.Range("J126").Value = 5977.143
$ws.Range("K126").Value = 5254.200000000001
$ws.Range("L126").Value = 17931.429
$ws.Range("M126").Value = -2784.200000000001
$ws.Range("N126").Value = -22871.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 22183.273
$ws.Range("J39").Value = 22183.273
$ws.Range("L39").Value = 66549.819
$ws.Range("N39").Value = -67137.819
$ws.Range("H68").Value = 2752.7324
$ws.Range("I68").Value = 969.3333
$ws.Range("J68").Value = 3501.76
$ws.Range("K68").Value = 2907.9999
$ws.Range("L68").Value = 10505.28
$ws.Range("M68").Value = -2096.9999
$ws.Range("N68").Value = -12127.28
$ws.Range("H71").Value = 2752.7324
$ws.Range("I71").Value = 969.3333
$ws.Range("J71").Value = 3501.76
$ws.Range("K71").Value = 8723.9997
$ws.Range("L71").Value = 31515.84
$ws.Range("M71").Value = -4667.9997
$ws.Range("N71").Value = -39627.84
$ws.Range("H81").Value = 1789.125
$ws.Range("J81").Value = 2700
$ws.Range("L81").Value = 8100
$ws.Range("N81").Value = -10346
$ws.Range("H84").Value = 1789.125
$ws.Range("J84").Value = 2700
$ws.Range("L84").Value = 24300
$ws.Range("N84").Value = -35532
$ws.Range("H116").Value = 8516
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 8516
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 25548
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -32432
$ws.Range("H122").Value = 2146.9
$ws.Range("I122").Value = 730.9474
$ws.Range("J122").Value = 2803.0732
$ws.Range("K122").Value = 6578.5266
$ws.Range("L122").Value = 25227.6588
$ws.Range("M122").Value = -4128.5266
$ws.Range("N122").Value = -30127.6588
$ws.Range("H131").Value = 864.1313
$ws.Range("I131").Value = 784
$ws.Range("J131").Value = 864.949
$ws.Range("K131").Value = 2352
$ws.Range("L131").Value = 2594.847
$ws.Range("M131").Value = 2688
$ws.Range("N131").Value = -12674.847

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 11111759
$ws.Range("I107").Value = 353.42856
$ws.Range("J107").Value = 37038372
$ws.Range("K107").Value = 353.42856
$ws.Range("L107").Value = 37038372
$ws.Range("M107").Value = 1566.57144
$ws.Range("N107").Value = -37042212
$ws.Range("H122").Value = 4570.4707
$ws.Range("J122").Value = 6198
$ws.Range("L122").Value = 18594
$ws.Range("N122").Value = -23494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5061.846
$ws.Range("I7").Value = 1879.8
$ws.Range("J7").Value = 7050.625
$ws.Range("K7").Value = 1879.8
$ws.Range("L7").Value = 7050.625
$ws.Range("M7").Value = -1767.8
$ws.Range("N7").Value = -7274.625
$ws.Range("H61").Value = 769.95
$ws.Range("I61").Value = 763.4706
$ws.Range("J61").Value = 806.6667
$ws.Range("K61").Value = 763.4706
$ws.Range("L61").Value = 806.6667
$ws.Range("M61").Value = -561.4706
$ws.Range("N61").Value = -1210.6667
$ws.Range("H113").Value = 769.95
$ws.Range("I113").Value = 763.4706
$ws.Range("J113").Value = 806.6667
$ws.Range("K113").Value = 763.4706
$ws.Range("L113").Value = 806.6667
$ws.Range("M113").Value = 1406.5294
$ws.Range("N113").Value = -5146.6667
$ws.Range("H122").Value = 5590.615
$ws.Range("I122").Value = 3222.25
$ws.Range("K122").Value = 9666.75
$ws.Range("M122").Value = -7216.75
$ws.Range("H126").Value = 5061.846
$ws.Range("I126").Value = 1879.8
$ws.Range("J126").Value = 7050.625
$ws.Range("K126").Value = 5639.4
$ws.Range("L126").Value = 21151.875
$ws.Range("M126").Value = -3169.4
$ws.Range("N126").Value = -26091.875
$ws.Range("H127").Value = 24746.666
$ws.Range("J127").Value = 24746.666
$ws.Range("L127").Value = 24746.666
$ws.Range("N127").Value = -34666.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3861.3333
$ws.Range("I122").Value = 2274
$ws.Range("J122").Value = 5131.2
$ws.Range("K122").Value = 6822
$ws.Range("L122").Value = 15393.6
$ws.Range("M122").Value = -4372
$ws.Range("N122").Value = -20293.6
